# Translate the row labels in column A from Chinese to English.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Data of China airline"
$ws.Range("A2").Value = "date"
$ws.Range("A3").Value = "Income"
$ws.Range("A4").Value = "Load factor"
$ws.Range("A5").Value = "Outcome"
$ws.Range("A6").Value = "Numbers of airplane(work)"
$ws.Range("A7").Value = "Number of patients"
$ws.Range("A8").Value = "New"

# Switch the workbook's default font away from the CJK "DengXian" to Calibri.
$wb.Styles("常规").Font.Name = "Calibri"

# Reset the lingering selection back to the top-left cell.
$ws.Range("A1").Select() | Out-Null
